# Segunda parte de implementacion del modulo administrar sesiones
# Set the value of Q12 on "Casos de Uso" sheet, which cascades through the
# shared formulas in that row (R12, U12, X12, ... AZ12, BA12).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

$ws.Range("Q12").Value = 2

# Update the active selection to reflect where the user last clicked (T12),
# matching the saved sheet view state.
$ws.Activate()
$ws.Range("T12").Select()
